$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Time Series image filename for the Bottom Temperature row (row 6)
$ws.Range("D6").Value = "BottomT_2025-04-17.png"

# Move the active selection to D6 (matches the updated cell)
$ws.Range("D6").Select()
